$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column D whose "Automatic Connection Feature Completed ..." text
# is reverted back to the plain "Automatic Connection Feature Completed".
$rows = @(3,4,6,8,9,11,19,20,21,22,23,24,25,27,28,31,34,35,37)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "Automatic Connection Feature Completed"
}

# Column D auto-fits narrower now that the long text is gone.
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).ColumnWidth = 35

# The active selection moved to J37 in the saved view.
$ws.Range("J37").Select() | Out-Null
